$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 12,20
$arr[0,0] = "ECs"
$arr[0,1] = "Agrn"
$arr[0,2] = "Musk"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 15.140316
$arr[0,7] = 45.420948
$arr[0,8] = 0.3229157245229468
$arr[0,9] = 0.3229157245229468
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 0.2865536666666667
$arr[0,13] = 0.859661
$arr[0,14] = 0.01543253808802733
$arr[0,15] = 0.01543253808802733
$arr[0,16] = 4.338513064291999
$arr[0,17] = 39.04661757862799
$arr[0,18] = 0.004983409217923319
$arr[0,19] = 0.004983409217923319
$arr[1,0] = "ECs"
$arr[1,1] = "Agrn"
$arr[1,2] = "Musk"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 15.140316
$arr[1,7] = 45.420948
$arr[1,8] = 0.3229157245229468
$arr[1,9] = 0.3229157245229468
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 9.389971666666668
$arr[1,13] = 28.169915
$arr[1,14] = 0.5057031622627903
$arr[1,15] = 0.5057031622627903
$arr[1,16] = 142.16713826438
$arr[1,17] = 1279.50424437942
$arr[1,18] = 0.1632995030356343
$arr[1,19] = 0.1632995030356343
$arr[2,0] = "ECs"
$arr[2,1] = "Agrn"
$arr[2,2] = "Musk"
$arr[2,3] = "sCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 15.140316
$arr[2,7] = 45.420948
$arr[2,8] = 0.3229157245229468
$arr[2,9] = 0.3229157245229468
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 8.891623666666666
$arr[2,13] = 26.674871
$arr[2,14] = 0.4788642996491824
$arr[2,15] = 0.4788642996491824
$arr[2,16] = 134.621992066412
$arr[2,17] = 1211.597928597708
$arr[2,18] = 0.1546328122693892
$arr[2,19] = 0.1546328122693892
$arr[3,0] = "FAPs"
$arr[3,1] = "Agrn"
$arr[3,2] = "Musk"
$arr[3,3] = "ECs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 18.94069966666667
$arr[3,7] = 56.822099
$arr[3,8] = 0.4039710766824948
$arr[3,9] = 0.4039710766824948
$arr[3,10] = 2
$arr[3,11] = 0.6666666666666666
$arr[3,12] = 0.2865536666666667
$arr[3,13] = 0.859661
$arr[3,14] = 0.01543253808802733
$arr[3,15] = 0.01543253808802733
$arr[3,16] = 5.427526938715444
$arr[3,17] = 48.847742448439
$arr[3,18] = 0.006234299027364013
$arr[3,19] = 0.006234299027364013
$arr[4,0] = "FAPs"
$arr[4,1] = "Agrn"
$arr[4,2] = "Musk"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 18.94069966666667
$arr[4,7] = 56.822099
$arr[4,8] = 0.4039710766824948
$arr[4,9] = 0.4039710766824948
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 9.389971666666668
$arr[4,13] = 28.169915
$arr[4,14] = 0.5057031622627903
$arr[4,15] = 0.5057031622627903
$arr[4,16] = 177.8526332168428
$arr[4,17] = 1600.673698951585
$arr[4,18] = 0.2042894509410418
$arr[4,19] = 0.2042894509410418
$arr[5,0] = "FAPs"
$arr[5,1] = "Agrn"
$arr[5,2] = "Musk"
$arr[5,3] = "sCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 18.94069966666667
$arr[5,7] = 56.822099
$arr[5,8] = 0.4039710766824948
$arr[5,9] = 0.4039710766824948
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 8.891623666666666
$arr[5,13] = 26.674871
$arr[5,14] = 0.4788642996491824
$arr[5,15] = 0.4788642996491824
$arr[5,16] = 168.4135734193588
$arr[5,17] = 1515.722160774229
$arr[5,18] = 0.193447326714089
$arr[5,19] = 0.193447326714089
$arr[6,0] = "M2"
$arr[6,1] = "Agrn"
$arr[6,2] = "Musk"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 3.221232
$arr[6,7] = 9.663696
$arr[6,8] = 0.06870308817441464
$arr[6,9] = 0.06870308817441464
$arr[6,10] = 2
$arr[6,11] = 0.6666666666666666
$arr[6,12] = 0.2865536666666667
$arr[6,13] = 0.859661
$arr[6,14] = 0.01543253808802733
$arr[6,15] = 0.01543253808802733
$arr[6,16] = 0.9230558407839999
$arr[6,17] = 8.307502567056
$arr[6,18] = 0.001060263025016754
$arr[6,19] = 0.001060263025016754
$arr[7,0] = "M2"
$arr[7,1] = "Agrn"
$arr[7,2] = "Musk"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 3.221232
$arr[7,7] = 9.663696
$arr[7,8] = 0.06870308817441464
$arr[7,9] = 0.06870308817441464
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 9.389971666666668
$arr[7,13] = 28.169915
$arr[7,14] = 0.5057031622627903
$arr[7,15] = 0.5057031622627903
$arr[7,16] = 30.24727721176
$arr[7,17] = 272.22549490584
$arr[7,18] = 0.0347433689470208
$arr[7,19] = 0.0347433689470208
$arr[8,0] = "M2"
$arr[8,1] = "Agrn"
$arr[8,2] = "Musk"
$arr[8,3] = "sCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 3.221232
$arr[8,7] = 9.663696
$arr[8,8] = 0.06870308817441464
$arr[8,9] = 0.06870308817441464
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 8.891623666666666
$arr[8,13] = 26.674871
$arr[8,14] = 0.4788642996491824
$arr[8,15] = 0.4788642996491824
$arr[8,16] = 28.641982687024
$arr[8,17] = 257.777844183216
$arr[8,18] = 0.03289945620237709
$arr[8,19] = 0.03289945620237709
$arr[9,0] = "sCs"
$arr[9,1] = "Agrn"
$arr[9,2] = "Musk"
$arr[9,3] = "ECs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 9.584029000000001
$arr[9,7] = 28.752087
$arr[9,8] = 0.2044101106201438
$arr[9,9] = 0.2044101106201438
$arr[9,10] = 2
$arr[9,11] = 0.6666666666666666
$arr[9,12] = 0.2865536666666667
$arr[9,13] = 0.859661
$arr[9,14] = 0.01543253808802733
$arr[9,15] = 0.01543253808802733
$arr[9,16] = 2.746338651389667
$arr[9,17] = 24.717047862507
$arr[9,18] = 0.00315456681772325
$arr[9,19] = 0.00315456681772325
$arr[10,0] = "sCs"
$arr[10,1] = "Agrn"
$arr[10,2] = "Musk"
$arr[10,3] = "FAPs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 9.584029000000001
$arr[10,7] = 28.752087
$arr[10,8] = 0.2044101106201438
$arr[10,9] = 0.2044101106201438
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 9.389971666666668
$arr[10,13] = 28.169915
$arr[10,14] = 0.5057031622627903
$arr[10,15] = 0.5057031622627903
$arr[10,16] = 89.99376076251168
$arr[10,17] = 809.9438468626051
$arr[10,18] = 0.1033708393390935
$arr[10,19] = 0.1033708393390935
$arr[11,0] = "sCs"
$arr[11,1] = "Agrn"
$arr[11,2] = "Musk"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 9.584029000000001
$arr[11,7] = 28.752087
$arr[11,8] = 0.2044101106201438
$arr[11,9] = 0.2044101106201438
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 8.891623666666666
$arr[11,13] = 26.674871
$arr[11,14] = 0.4788642996491824
$arr[11,15] = 0.4788642996491824
$arr[11,16] = 85.21757907841967
$arr[11,17] = 766.9582117057771
$arr[11,18] = 0.09788470446332707
$arr[11,19] = 0.09788470446332707

$ws.Range("A2:T13").Value = $arr
